# Updates Leve profit-tracking data across the Ixion_Profits sheets
# (currentAveragePrice / LevePrice / LeveProfit columns), refreshed by the
# scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 600.3333
$ws.Range("I55").Value = 737.1111
$ws.Range("K55").Value = 737.1111
$ws.Range("M55").Value = -523.1111

$ws.Range("H137").Value = 1514.0426
$ws.Range("I137").Value = 1115.5238
$ws.Range("J137").Value = 1835.9231
$ws.Range("K137").Value = 3346.5714
$ws.Range("L137").Value = 5507.7693
$ws.Range("M137").Value = -796.5713999999998
$ws.Range("N137").Value = -10607.7693

$ws.Range("H141").Value = 1949.174
$ws.Range("I141").Value = 2044.9048
$ws.Range("J141").Value = 944
$ws.Range("K141").Value = 6134.7144
$ws.Range("L141").Value = 2832
$ws.Range("M141").Value = -954.7143999999998
$ws.Range("N141").Value = -13192

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4943.3945
$ws.Range("I32").Value = 4724.0713
$ws.Range("J32").Value = 5557.5
$ws.Range("K32").Value = 4724.0713
$ws.Range("L32").Value = 5557.5
$ws.Range("M32").Value = -4437.0713
$ws.Range("N32").Value = -6131.5

$ws.Range("H74").Value = 2554
$ws.Range("I74").Value = 2776.0981
$ws.Range("J74").Value = 1887.7059
$ws.Range("K74").Value = 2776.0981
$ws.Range("L74").Value = 1887.7059
$ws.Range("M74").Value = -1902.0981
$ws.Range("N74").Value = -3635.7059

$ws.Range("H77").Value = 2554
$ws.Range("I77").Value = 2776.0981
$ws.Range("J77").Value = 1887.7059
$ws.Range("K77").Value = 13880.4905
$ws.Range("L77").Value = 9438.529500000001
$ws.Range("M77").Value = -9512.4905
$ws.Range("N77").Value = -18174.5295

$ws.Range("H114").Value = 166666.33
$ws.Range("J114").Value = 166666.33
$ws.Range("L114").Value = 166666.33
$ws.Range("N114").Value = -175344.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 500003100
$ws.Range("J25").Value = 500003100
$ws.Range("L25").Value = 500003100
$ws.Range("N25").Value = -500003448

$ws.Range("H58").Value = 1636.4062
$ws.Range("I58").Value = 1092.55
$ws.Range("J58").Value = 2542.8333
$ws.Range("K58").Value = 1092.55
$ws.Range("L58").Value = 2542.8333
$ws.Range("M58").Value = -889.55
$ws.Range("N58").Value = -2948.8333

$ws.Range("H136").Value = 1636.4062
$ws.Range("I136").Value = 1092.55
$ws.Range("J136").Value = 2542.8333
$ws.Range("K136").Value = 3277.65
$ws.Range("L136").Value = 7628.499899999999
$ws.Range("M136").Value = -727.6499999999996
$ws.Range("N136").Value = -12728.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 161471.34
$ws.Range("I5").Value = 14719.571
$ws.Range("J5").Value = 194608.84
$ws.Range("K5").Value = 44158.713
$ws.Range("L5").Value = 583826.52
$ws.Range("M5").Value = -44046.713
$ws.Range("N5").Value = -584050.52

$ws.Range("H68").Value = 2750.8718
$ws.Range("I68").Value = 3564.3784
$ws.Range("J68").Value = 2016.7317
$ws.Range("K68").Value = 10693.1352
$ws.Range("L68").Value = 6050.1951
$ws.Range("M68").Value = -9882.135200000001
$ws.Range("N68").Value = -7672.1951

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H70").Value = 2750.25
$ws.Range("I70").Value = 2167
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 6501
$ws.Range("L70").Value = 13500
$ws.Range("M70").Value = -6186
$ws.Range("N70").Value = -14130

$ws.Range("H71").Value = 2750.8718
$ws.Range("I71").Value = 3564.3784
$ws.Range("J71").Value = 2016.7317
$ws.Range("K71").Value = 32079.4056
$ws.Range("L71").Value = 18150.5853
$ws.Range("M71").Value = -28023.4056
$ws.Range("N71").Value = -26262.5853

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H73").Value = 2750.25
$ws.Range("I73").Value = 2167
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 6501
$ws.Range("L73").Value = 13500
$ws.Range("M73").Value = -5409
$ws.Range("N73").Value = -15684

$ws.Range("H131").Value = 12223447
$ws.Range("I131").Value = 11111729
$ws.Range("J131").Value = 12346971
$ws.Range("K131").Value = 33335187
$ws.Range("L131").Value = 37040913
$ws.Range("M131").Value = -33330147
$ws.Range("N131").Value = -37050993

$ws.Range("H135").Value = 161471.34
$ws.Range("I135").Value = 14719.571
$ws.Range("J135").Value = 194608.84
$ws.Range("K135").Value = 132476.139
$ws.Range("L135").Value = 1751479.56
$ws.Range("M135").Value = -129941.139
$ws.Range("N135").Value = -1756549.56

$ws.Range("H137").Value = 30318538
$ws.Range("I137").Value = 1839.625
$ws.Range("J137").Value = 40019880
$ws.Range("K137").Value = 5518.875
$ws.Range("L137").Value = 120059640
$ws.Range("M137").Value = -418.875
$ws.Range("N137").Value = -120069840

$ws.Range("H139").Value = 4043.1135
$ws.Range("I139").Value = 4969.88
$ws.Range("J139").Value = 2823.6843
$ws.Range("K139").Value = 14909.64
$ws.Range("L139").Value = 8471.052899999999
$ws.Range("M139").Value = -9769.639999999999
$ws.Range("N139").Value = -18751.0529

$ws.Range("H140").Value = 3493.8333
$ws.Range("I140").Value = 3493.8333
$ws.Range("K140").Value = 10481.4999
$ws.Range("M140").Value = -5301.499899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 393.53333
$ws.Range("I2").Value = 401.64285
$ws.Range("J2").Value = 280
$ws.Range("K2").Value = 401.64285
$ws.Range("L2").Value = 280
$ws.Range("M2").Value = -288.64285
$ws.Range("N2").Value = -506

$ws.Range("H33").Value = 5800
$ws.Range("J33").Value = 5800
$ws.Range("L33").Value = 5800
$ws.Range("N33").Value = -6304

$ws.Range("H103").Value = 29000
$ws.Range("J103").Value = 29000
$ws.Range("L103").Value = 29000
$ws.Range("N103").Value = -31344

$ws.Range("H132").Value = 3564.7144
$ws.Range("I132").Value = 3795.875
$ws.Range("K132").Value = 11387.625
$ws.Range("M132").Value = -8857.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5558135.5
$ws.Range("I22").Value = 37038104
$ws.Range("J22").Value = 2847.0588
$ws.Range("K22").Value = 37038104
$ws.Range("L22").Value = 2847.0588
$ws.Range("M22").Value = -37037809
$ws.Range("N22").Value = -3437.0588

$ws.Range("H27").Value = 5558135.5
$ws.Range("I27").Value = 37038104
$ws.Range("J27").Value = 2847.0588
$ws.Range("K27").Value = 37038104
$ws.Range("L27").Value = 2847.0588
$ws.Range("M27").Value = -37037997
$ws.Range("N27").Value = -3061.0588

$ws.Range("H132").Value = 16674061
$ws.Range("I132").Value = 25499634
$ws.Range("J132").Value = 3534.7778
$ws.Range("K132").Value = 76498902
$ws.Range("L132").Value = 10604.3334
$ws.Range("M132").Value = -76496372
$ws.Range("N132").Value = -15664.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()

$ws.Range("H139").Value = 62000
$ws.Range("J139").Value = 62000
$ws.Range("L139").Value = 62000
$ws.Range("N139").Value = -72280
